$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 4 with the new "all-mpnet-base-v2" XGBoost model entry,
# overwriting the old filler/duplicate row (which used to just repeat the
# MiniLM model's data with a "subsample: 0.9" hyper-param placeholder).
# Fill order matches how the shared-string table grows (A, E, then B).
$ws.Range("A4").Value = "XGBoost + ""all-mpnet-base-v2"" BERT-Sentence Embedding"
$ws.Range("E4").Value = "lowercase,`npunctuation removal`nextra features`nno stemming"
$ws.Range("B4").Value = "'colsample_bytree': 0.8, 'learning_rate': np.float64(0.03686562370169114), 'max_depth': 7, 'min_child_weight': 3, 'n_estimators': 200, 'subsample': 0.7"
$ws.Range("C4").Value = 0.30713916629328603
$ws.Range("D4").Value = 0.86075589304707001

# The new row is a bit taller than the default 60pt used by the other rows.
$ws.Rows("4").RowHeight = 64.5

# Move the active selection to B3, matching the author's last edit position.
$ws.Range("B3").Select()

# Touch the page setup so a pageSetup element gets written for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
